$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> subject name (mapel) for the new full list
$rowSubjects = @{
    2  = "PAI"
    3  = "PPKN"
    4  = "Bahasa Indonesia"
    5  = "Matematika Wajib"
    6  = "Sejarah Indonesia"
    7  = "Bahasa Inggris"
    8  = "Seni Budaya"
    9  = "Penjasorkes"
    10 = "Prakarya"
    11 = "Bahasa Daerah"
    12 = "Matematika Peminatan"
    13 = "Biologi"
    14 = "Fisika"
    15 = "Kimia"
    16 = "Sejarah Peminatan"
    17 = "Geografi"
    18 = "Ekonomi"
    19 = "Sosiologi"
    20 = "Sastra Inggris"
}

# Order in which the rows were filled in (row 6 entered last)
$fillOrder = @(2, 3, 4, 5, 7, 8, 9, 10, 11, 16, 12, 13, 14, 15, 17, 18, 19, 20, 6)

foreach ($row in $fillOrder) {
    $ws.Cells.Item($row, 1).Value = $row - 1
    $ws.Cells.Item($row, 2).Value = $rowSubjects[$row]
}

$ws.Range("B9").Select()
